$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet to "Icons"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Icons"

# Add a new sheet right after "Icons" and name it "Sheet2"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Fill in the new "High Tier Uniques" label on the Icons sheet
$ws1.Range("C24").Value = "High Tier Uniques"

# Move the view / selection on Icons to match the saved state
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws1.Range("C24").Select()
